# Adds the "28. 9. 2021" wave of data to both sheets, matching the upstream
# ZBP_03_strategie_domacnosti update (new column AH on "data", AG on "pocetR").

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("data")
$ws2 = $wb.Worksheets.Item("pocetR")

# --- Sheet "data": new header date in AH1 (copy formatting from AG1, the
# previous last header cell, then set the new label) ---
$ws1.Range("AG1").Copy()
$ws1.Range("AH1").PasteSpecial(-4122)
$ws1.Range("AH1").Value = "28. 9. 2021"

# --- Sheet "data": new percentage values for rows 2-45 ---
$ws1.Range("AH2").Value = 0.18
$ws1.Range("AH3").Value = 0.11
$ws1.Range("AH4").Value = 0.41
$ws1.Range("AH5").Value = 0.28
$ws1.Range("AH6").Value = 0.13
$ws1.Range("AH7").Value = 0.17
$ws1.Range("AH8").Value = 0.19
$ws1.Range("AH9").Value = 0.19
$ws1.Range("AH10").Value = 0.2
$ws1.Range("AH11").Value = 0.17
$ws1.Range("AH12").Value = 0.18
$ws1.Range("AH13").Value = 0.37
$ws1.Range("AH14").Value = 0.18
$ws1.Range("AH15").Value = 0.17
$ws1.Range("AH16").Value = 0.21
$ws1.Range("AH17").Value = 0.19
$ws1.Range("AH18").Value = 0.18
$ws1.Range("AH19").Value = 0.24
$ws1.Range("AH20").Value = 0.16
$ws1.Range("AH21").Value = 0.14
$ws1.Range("AH22").Value = 0.12
$ws1.Range("AH23").Value = 0.21
$ws1.Range("AH24").Value = 0.37
$ws1.Range("AH25").Value = 0.39
$ws1.Range("AH26").Value = 0.13
$ws1.Range("AH27").Value = 0.07000000000000001
$ws1.Range("AH28").Value = 0.13
$ws1.Range("AH29").Value = 0.19
$ws1.Range("AH30").Value = 0.07000000000000001
$ws1.Range("AH31").Value = 0.09
$ws1.Range("AH32").Value = 0.13
$ws1.Range("AH33").Value = 0.17
$ws1.Range("AH34").Value = 0.23
$ws1.Range("AH35").Value = 0.06
$ws1.Range("AH36").Value = 0.11
$ws1.Range("AH37").Value = 0.13
$ws1.Range("AH38").Value = 0.07000000000000001
$ws1.Range("AH39").Value = 0.17
$ws1.Range("AH40").Value = 0.13
$ws1.Range("AH41").Value = 0.08
$ws1.Range("AH42").Value = 0.1
$ws1.Range("AH43").Value = 0.09
$ws1.Range("AH44").Value = 0.19
$ws1.Range("AH45").Value = 0.25

# --- Sheet "data": update the "aktualizace" footnote in row 46 ---
$ws1.Range("A46").Value = "Život během pandemie, Strategie domácností, % respondentů celkově a ve skupinách, aktualizace 6. 10. 2021"

# --- Sheet "pocetR": new header date in AG1 (copy formatting from AF1) ---
$ws2.Range("AF1").Copy()
$ws2.Range("AG1").PasteSpecial(-4122)
$ws2.Range("AG1").Value = "28. 9. 2021"

# --- Sheet "pocetR": new sample-size values for rows 2-23 ---
$ws2.Range("AG2").Value = 1855
$ws2.Range("AG3").Value = 165
$ws2.Range("AG4").Value = 370
$ws2.Range("AG5").Value = 1320
$ws2.Range("AG6").Value = 888
$ws2.Range("AG7").Value = 165
$ws2.Range("AG8").Value = 531
$ws2.Range("AG9").Value = 271
$ws2.Range("AG10").Value = 853
$ws2.Range("AG11").Value = 154
$ws2.Range("AG12").Value = 105
$ws2.Range("AG13").Value = 743
$ws2.Range("AG14").Value = 859
$ws2.Range("AG15").Value = 631
$ws2.Range("AG16").Value = 365
$ws2.Range("AG17").Value = 180
$ws2.Range("AG18").Value = 679
$ws2.Range("AG19").Value = 613
$ws2.Range("AG20").Value = 255
$ws2.Range("AG21").Value = 566
$ws2.Range("AG22").Value = 309
$ws2.Range("AG23").Value = 168

# --- Sheet "pocetR": keep row 24 fully populated through column AG (blank
# placeholder cell, same as B24:AF24) and update its footnote text ---
$ws2.Range("AG24").Font.Bold = $false
$ws2.Range("A24").Value = "Život během pandemie, Strategie domácností, velikost dotázaného souboru celkově a ve skupinách, aktualizace 6. 10. 2021"
